# Regenerate the "K" column (G) values on Sheet1.
# The author's pipeline recomputed the K statistic (previously "Strike#")
# and rewrote the save_data sheet; here we apply the resulting cell values
# directly to column G for each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 7
    4  = 0
    5  = 1
    6  = 1
    7  = 0
    8  = 1
    9  = 0
    10 = 1
    11 = 2
    12 = 4
    13 = 2
    14 = 2
    15 = 5
    16 = 4
    17 = 0
    18 = 0
    19 = 1
    20 = 2
    21 = 1
    22 = 2
    23 = 2
    24 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
